$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new member row (row 3) with the author's information.
# Columns: A=Full Name, B=Github, C=Facebook, D=Role, E=Team
$ws.Range("A3").Value = "Samir"
$ws.Range("B3").Value = "Samir-SB"
$ws.Range("C3").Value = "Samir-SB"
$ws.Range("D3").Value = "member"
$ws.Range("E3").Value = "git team"

# Update the selection to reflect where the user left off editing.
$ws.Range("D4").Select()
